$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 14: "Fetch 'initial load' of open Slots* for Service(s)" diagram.
# Ungroup "Group 14", drop the small step-label textbox ("AA"), leave the
# connector + textbox where the ungroup operation naturally places them.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$grp14 = $s14.Shapes.Item("Group 14")
$members14 = $grp14.Ungroup()
for ($i = 1; $i -le $members14.Count; $i++) {
    $sh = $members14.Item($i)
    if ($sh.Name -eq "TextBox 28") {
        $sh.Delete()
    }
}

# ---------------------------------------------------------------------------
# Slide 15: "Fetch all Slots* updated since last fetch" diagram.
# Ungroup "Group 14", drop the small step-label textbox ("C"), and move /
# resize the remaining description textbox to its final authored position.
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$grp15 = $s15.Shapes.Item("Group 14")
$members15 = $grp15.Ungroup()
for ($i = 1; $i -le $members15.Count; $i++) {
    $sh = $members15.Item($i)
    if ($sh.Name -eq "TextBox 28") {
        $sh.Delete()
    } elseif ($sh.Name -eq "TextBox 3") {
        $sh.Left = 516.8725590551
        $sh.Top = 207.6953937008
        $sh.Width = 199.7740551181
        $sh.Height = 50.8922440945
    }
}

# ---------------------------------------------------------------------------
# Slide 8: "Share Business Rules" diagram.
# Ungroup "Group 14", drop the small step-label textbox ("A"), and move the
# remaining description textbox to its final authored position (size keeps
# the value produced by the ungroup transform).
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$grp8 = $s8.Shapes.Item("Group 14")
$members8 = $grp8.Ungroup()
for ($i = 1; $i -le $members8.Count; $i++) {
    $sh = $members8.Item($i)
    if ($sh.Name -eq "TextBox 28") {
        $sh.Delete()
    } elseif ($sh.Name -eq "TextBox 3") {
        $sh.Left = 536.3202755906
        $sh.Top = 243.7694881890
    }
}
